$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Motor row (row 3) with the new gearmotor variant details.
$ws.Range("A3").Value = "99:1 Metal Gearmotor 25Dx69L mm HP 12V with 48 CPR Encoder"
$ws.Range("C3").Value = 4847
$ws.Range("D3").Value = 58
$ws.Range("F3").Value = "SKU-005518"
$ws.Range("G3").Value = "100RPM, 2.11 Nm, 12V 5A"

# Update the selected/active cell to A3.
$ws.Range("A3").Select()
